# Auto-generated edit script applying the 'Horarios actualizados Linea 141 - 605' update
# Updates header info rows (2,3) and the scraped schedule rows for each of the 3 sheets.

$wb = $excel.ActiveWorkbook

# ----- Sheet: LP1912 -----
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 1).Value = "Última actualización: 17:55:25"

$ws.Cells.Item(3, 1).Value = "Total filas: 301"

$ws.Cells.Item(118, 1).Value = "11:33:52"
$ws.Cells.Item(118, 2).Value = "11:52"
$ws.Cells.Item(118, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(118, 4).Value = 19
$ws.Cells.Item(118, 5).Value = "LP1912"

$ws.Cells.Item(119, 1).Value = "10:36:50"
$ws.Cells.Item(119, 2).Value = "11:52"
$ws.Cells.Item(119, 3).Value = "225_GOMEZ"
$ws.Cells.Item(119, 4).Value = 76
$ws.Cells.Item(119, 5).Value = "LP1912"

$ws.Cells.Item(137, 1).Value = "10:36:50"
$ws.Cells.Item(137, 2).Value = "12:34"
$ws.Cells.Item(137, 3).Value = "15_ABASTO"
$ws.Cells.Item(137, 4).Value = 118
$ws.Cells.Item(137, 5).Value = "LP1912"

$ws.Cells.Item(138, 1).Value = "11:46:32"
$ws.Cells.Item(138, 2).Value = "12:34"
$ws.Cells.Item(138, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(138, 4).Value = 48
$ws.Cells.Item(138, 5).Value = "LP1912"

$ws.Cells.Item(147, 1).Value = "11:13:15"
$ws.Cells.Item(147, 2).Value = "13:03"
$ws.Cells.Item(147, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(147, 4).Value = 110
$ws.Cells.Item(147, 5).Value = "LP1912"

$ws.Cells.Item(148, 1).Value = "11:33:52"
$ws.Cells.Item(148, 2).Value = "13:03"
$ws.Cells.Item(148, 3).Value = "215C_EL PATO"
$ws.Cells.Item(148, 4).Value = 90
$ws.Cells.Item(148, 5).Value = "LP1912"

$ws.Cells.Item(220, 1).Value = "15:16:46"
$ws.Cells.Item(220, 2).Value = "16:30"
$ws.Cells.Item(220, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(220, 4).Value = 74
$ws.Cells.Item(220, 5).Value = "LP1912"

$ws.Cells.Item(221, 1).Value = "16:12:06"
$ws.Cells.Item(221, 2).Value = "16:30"
$ws.Cells.Item(221, 3).Value = "14_ABASTO"
$ws.Cells.Item(221, 4).Value = 18
$ws.Cells.Item(221, 5).Value = "LP1912"

$ws.Cells.Item(260, 1).Value = "17:55:25"
$ws.Cells.Item(260, 2).Value = "17:55"
$ws.Cells.Item(260, 3).Value = "10_OLMOS"
$ws.Cells.Item(260, 4).Value = 0
$ws.Cells.Item(260, 5).Value = "LP1912"

$ws.Cells.Item(261, 1).Value = "16:44:58"
$ws.Cells.Item(261, 2).Value = "17:57"
$ws.Cells.Item(261, 3).Value = "17_ROMERO"
$ws.Cells.Item(261, 4).Value = 73
$ws.Cells.Item(261, 5).Value = "LP1912"

$ws.Cells.Item(262, 1).Value = "16:12:06"
$ws.Cells.Item(262, 2).Value = "17:58"
$ws.Cells.Item(262, 3).Value = "17_ROMERO"
$ws.Cells.Item(262, 4).Value = 106
$ws.Cells.Item(262, 5).Value = "LP1912"

$ws.Cells.Item(263, 1).Value = "16:28:21"
$ws.Cells.Item(263, 2).Value = "18:05"
$ws.Cells.Item(263, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(263, 4).Value = 97
$ws.Cells.Item(263, 5).Value = "LP1912"

$ws.Cells.Item(264, 1).Value = "16:12:06"
$ws.Cells.Item(264, 2).Value = "18:06"
$ws.Cells.Item(264, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(264, 4).Value = 114
$ws.Cells.Item(264, 5).Value = "LP1912"

$ws.Cells.Item(265, 1).Value = "16:44:58"
$ws.Cells.Item(265, 2).Value = "18:09"
$ws.Cells.Item(265, 3).Value = "15_ABASTO"
$ws.Cells.Item(265, 4).Value = 85
$ws.Cells.Item(265, 5).Value = "LP1912"

$ws.Cells.Item(266, 1).Value = "16:44:58"
$ws.Cells.Item(266, 2).Value = "18:09"
$ws.Cells.Item(266, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(266, 4).Value = 85
$ws.Cells.Item(266, 5).Value = "LP1912"

$ws.Cells.Item(267, 1).Value = "16:12:06"
$ws.Cells.Item(267, 2).Value = "18:10"
$ws.Cells.Item(267, 3).Value = "15_ABASTO"
$ws.Cells.Item(267, 4).Value = 118
$ws.Cells.Item(267, 5).Value = "LP1912"

$ws.Cells.Item(268, 1).Value = "16:12:06"
$ws.Cells.Item(268, 2).Value = "18:10"
$ws.Cells.Item(268, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(268, 4).Value = 118
$ws.Cells.Item(268, 5).Value = "LP1912"

$ws.Cells.Item(269, 1).Value = "16:44:58"
$ws.Cells.Item(269, 2).Value = "18:16"
$ws.Cells.Item(269, 3).Value = "10_OLMOS"
$ws.Cells.Item(269, 4).Value = 92
$ws.Cells.Item(269, 5).Value = "LP1912"

$ws.Cells.Item(270, 1).Value = "16:28:21"
$ws.Cells.Item(270, 2).Value = "18:17"
$ws.Cells.Item(270, 3).Value = "10_OLMOS"
$ws.Cells.Item(270, 4).Value = 109
$ws.Cells.Item(270, 5).Value = "LP1912"

$ws.Cells.Item(271, 1).Value = "16:37:37"
$ws.Cells.Item(271, 2).Value = "18:21"
$ws.Cells.Item(271, 3).Value = "215C_EL PATO"
$ws.Cells.Item(271, 4).Value = 104
$ws.Cells.Item(271, 5).Value = "LP1912"

$ws.Cells.Item(272, 1).Value = "16:28:21"
$ws.Cells.Item(272, 2).Value = "18:22"
$ws.Cells.Item(272, 3).Value = "215C_EL PATO"
$ws.Cells.Item(272, 4).Value = 114
$ws.Cells.Item(272, 5).Value = "LP1912"

$ws.Cells.Item(273, 1).Value = "16:28:21"
$ws.Cells.Item(273, 2).Value = "18:25"
$ws.Cells.Item(273, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(273, 4).Value = 117
$ws.Cells.Item(273, 5).Value = "LP1912"

$ws.Cells.Item(274, 1).Value = "17:13:30"
$ws.Cells.Item(274, 2).Value = "18:29"
$ws.Cells.Item(274, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(274, 4).Value = 76
$ws.Cells.Item(274, 5).Value = "LP1912"

$ws.Cells.Item(275, 1).Value = "16:37:37"
$ws.Cells.Item(275, 2).Value = "18:29"
$ws.Cells.Item(275, 3).Value = "14_ABASTO"
$ws.Cells.Item(275, 4).Value = 112
$ws.Cells.Item(275, 5).Value = "LP1912"

$ws.Cells.Item(276, 1).Value = "17:55:25"
$ws.Cells.Item(276, 2).Value = "18:30"
$ws.Cells.Item(276, 3).Value = "14_ABASTO"
$ws.Cells.Item(276, 4).Value = 35
$ws.Cells.Item(276, 5).Value = "LP1912"

$ws.Cells.Item(277, 1).Value = "17:47:45"
$ws.Cells.Item(277, 2).Value = "18:34"
$ws.Cells.Item(277, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(277, 4).Value = 47
$ws.Cells.Item(277, 5).Value = "LP1912"

$ws.Cells.Item(278, 1).Value = "16:44:58"
$ws.Cells.Item(278, 2).Value = "18:35"
$ws.Cells.Item(278, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(278, 4).Value = 111
$ws.Cells.Item(278, 5).Value = "LP1912"

$ws.Cells.Item(279, 1).Value = "16:37:37"
$ws.Cells.Item(279, 2).Value = "18:36"
$ws.Cells.Item(279, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(279, 4).Value = 119
$ws.Cells.Item(279, 5).Value = "LP1912"

$ws.Cells.Item(280, 1).Value = "17:35:41"
$ws.Cells.Item(280, 2).Value = "18:37"
$ws.Cells.Item(280, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(280, 4).Value = 62
$ws.Cells.Item(280, 5).Value = "LP1912"

$ws.Cells.Item(281, 1).Value = "16:44:58"
$ws.Cells.Item(281, 2).Value = "18:40"
$ws.Cells.Item(281, 3).Value = "10_OLMOS"
$ws.Cells.Item(281, 4).Value = 116
$ws.Cells.Item(281, 5).Value = "LP1912"

$ws.Cells.Item(282, 1).Value = "17:13:30"
$ws.Cells.Item(282, 2).Value = "18:41"
$ws.Cells.Item(282, 3).Value = "10_OLMOS"
$ws.Cells.Item(282, 4).Value = 88
$ws.Cells.Item(282, 5).Value = "LP1912"

$ws.Cells.Item(283, 1).Value = "17:13:30"
$ws.Cells.Item(283, 2).Value = "18:45"
$ws.Cells.Item(283, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(283, 4).Value = 92
$ws.Cells.Item(283, 5).Value = "LP1912"

$ws.Cells.Item(284, 1).Value = "17:55:25"
$ws.Cells.Item(284, 2).Value = "18:46"
$ws.Cells.Item(284, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(284, 4).Value = 51
$ws.Cells.Item(284, 5).Value = "LP1912"

$ws.Cells.Item(285, 1).Value = "17:13:30"
$ws.Cells.Item(285, 2).Value = "18:52"
$ws.Cells.Item(285, 3).Value = "17_ROMERO"
$ws.Cells.Item(285, 4).Value = 99
$ws.Cells.Item(285, 5).Value = "LP1912"

$ws.Cells.Item(286, 1).Value = "17:13:30"
$ws.Cells.Item(286, 2).Value = "18:57"
$ws.Cells.Item(286, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(286, 4).Value = 104
$ws.Cells.Item(286, 5).Value = "LP1912"

$ws.Cells.Item(287, 1).Value = "17:13:30"
$ws.Cells.Item(287, 2).Value = "18:59"
$ws.Cells.Item(287, 3).Value = "14_ABASTO"
$ws.Cells.Item(287, 4).Value = 106
$ws.Cells.Item(287, 5).Value = "LP1912"

$ws.Cells.Item(288, 1).Value = "17:47:45"
$ws.Cells.Item(288, 2).Value = "19:02"
$ws.Cells.Item(288, 3).Value = "14_ABASTO"
$ws.Cells.Item(288, 4).Value = 75
$ws.Cells.Item(288, 5).Value = "LP1912"

$ws.Cells.Item(289, 1).Value = "17:35:41"
$ws.Cells.Item(289, 2).Value = "19:03"
$ws.Cells.Item(289, 3).Value = "215_EL PELIGRO"
$ws.Cells.Item(289, 4).Value = 88
$ws.Cells.Item(289, 5).Value = "LP1912"

$ws.Cells.Item(290, 1).Value = "17:55:25"
$ws.Cells.Item(290, 2).Value = "19:03"
$ws.Cells.Item(290, 3).Value = "14_ABASTO"
$ws.Cells.Item(290, 4).Value = 68
$ws.Cells.Item(290, 5).Value = "LP1912"

$ws.Cells.Item(291, 1).Value = "17:13:30"
$ws.Cells.Item(291, 2).Value = "19:04"
$ws.Cells.Item(291, 3).Value = "215_EL PELIGRO"
$ws.Cells.Item(291, 4).Value = 111
$ws.Cells.Item(291, 5).Value = "LP1912"

$ws.Cells.Item(292, 1).Value = "17:55:25"
$ws.Cells.Item(292, 2).Value = "19:11"
$ws.Cells.Item(292, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(292, 4).Value = 76
$ws.Cells.Item(292, 5).Value = "LP1912"

$ws.Cells.Item(293, 1).Value = "17:55:25"
$ws.Cells.Item(293, 2).Value = "19:14"
$ws.Cells.Item(293, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(293, 4).Value = 79
$ws.Cells.Item(293, 5).Value = "LP1912"

$ws.Cells.Item(294, 1).Value = "17:47:45"
$ws.Cells.Item(294, 2).Value = "19:15"
$ws.Cells.Item(294, 3).Value = "17_ROMERO"
$ws.Cells.Item(294, 4).Value = 88
$ws.Cells.Item(294, 5).Value = "LP1912"

$ws.Cells.Item(295, 1).Value = "17:35:41"
$ws.Cells.Item(295, 2).Value = "19:16"
$ws.Cells.Item(295, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(295, 4).Value = 101
$ws.Cells.Item(295, 5).Value = "LP1912"

$ws.Cells.Item(296, 1).Value = "17:55:25"
$ws.Cells.Item(296, 2).Value = "19:16"
$ws.Cells.Item(296, 3).Value = "17_ROMERO"
$ws.Cells.Item(296, 4).Value = 81
$ws.Cells.Item(296, 5).Value = "LP1912"

$ws.Cells.Item(297, 1).Value = "17:35:41"
$ws.Cells.Item(297, 2).Value = "19:17"
$ws.Cells.Item(297, 3).Value = "14X44_ABASTO"
$ws.Cells.Item(297, 4).Value = 102
$ws.Cells.Item(297, 5).Value = "LP1912"

$ws.Cells.Item(298, 1).Value = "17:55:25"
$ws.Cells.Item(298, 2).Value = "19:22"
$ws.Cells.Item(298, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(298, 4).Value = 87
$ws.Cells.Item(298, 5).Value = "LP1912"

$ws.Cells.Item(299, 1).Value = "17:35:41"
$ws.Cells.Item(299, 2).Value = "19:27"
$ws.Cells.Item(299, 3).Value = "215C_EL PATO"
$ws.Cells.Item(299, 4).Value = 112
$ws.Cells.Item(299, 5).Value = "LP1912"

$ws.Cells.Item(300, 1).Value = "17:55:25"
$ws.Cells.Item(300, 2).Value = "19:28"
$ws.Cells.Item(300, 3).Value = "215C_EL PATO"
$ws.Cells.Item(300, 4).Value = 93
$ws.Cells.Item(300, 5).Value = "LP1912"

$ws.Cells.Item(301, 1).Value = "17:47:45"
$ws.Cells.Item(301, 2).Value = "19:35"
$ws.Cells.Item(301, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(301, 4).Value = 108
$ws.Cells.Item(301, 5).Value = "LP1912"

$ws.Cells.Item(302, 1).Value = "17:55:25"
$ws.Cells.Item(302, 2).Value = "19:36"
$ws.Cells.Item(302, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(302, 4).Value = 101
$ws.Cells.Item(302, 5).Value = "LP1912"

$ws.Cells.Item(303, 1).Value = "17:55:25"
$ws.Cells.Item(303, 2).Value = "19:39"
$ws.Cells.Item(303, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(303, 4).Value = 104
$ws.Cells.Item(303, 5).Value = "LP1912"

$ws.Cells.Item(304, 1).Value = "17:47:45"
$ws.Cells.Item(304, 2).Value = "19:42"
$ws.Cells.Item(304, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(304, 4).Value = 115
$ws.Cells.Item(304, 5).Value = "LP1912"

$ws.Cells.Item(305, 1).Value = "17:55:25"
$ws.Cells.Item(305, 2).Value = "19:52"
$ws.Cells.Item(305, 3).Value = "81_EL PELIGRO"
$ws.Cells.Item(305, 4).Value = 117
$ws.Cells.Item(305, 5).Value = "LP1912"

$ws.Cells.Item(306, 1).Value = "17:55:25"
$ws.Cells.Item(306, 2).Value = "19:53"
$ws.Cells.Item(306, 3).Value = "225_GOMEZ"
$ws.Cells.Item(306, 4).Value = 118
$ws.Cells.Item(306, 5).Value = "LP1912"

# ----- Sheet: LP1912-215 -----
$ws = $wb.Worksheets.Item(2)

$ws.Cells.Item(2, 1).Value = "Última actualización: 17:55:25"

$ws.Cells.Item(3, 1).Value = "Total filas: 50"

$ws.Cells.Item(55, 1).Value = "17:55:25"
$ws.Cells.Item(55, 2).Value = "19:28"
$ws.Cells.Item(55, 3).Value = "215C_EL PATO"
$ws.Cells.Item(55, 4).Value = 93
$ws.Cells.Item(55, 5).Value = "LP1912"

# ----- Sheet: 6203-6173 -----
$ws = $wb.Worksheets.Item(3)

$ws.Cells.Item(2, 1).Value = "Última actualización: 17:55:25"

$ws.Cells.Item(3, 1).Value = "Total filas: 43"

$ws.Cells.Item(48, 1).Value = "17:55:25"
$ws.Cells.Item(48, 2).Value = "19:24"
$ws.Cells.Item(48, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(48, 4).Value = 89
$ws.Cells.Item(48, 5).Value = "L6173"
